# Auto-generated edit script applying the numeric corrections described in the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 68629064
$ws.Range("I43").Value = 166667800
$ws.Range("J43").Value = 15153385
$ws.Range("K43").Value = 166667800
$ws.Range("L43").Value = 15153385
$ws.Range("M43").Value = -166667731
$ws.Range("N43").Value = -15153523
$ws.Range("H112").Value = 1536.1538
$ws.Range("I112").Value = 650
$ws.Range("J112").Value = 1697.2727
$ws.Range("K112").Value = 1950
$ws.Range("L112").Value = 5091.8181
$ws.Range("M112").Value = -842
$ws.Range("N112").Value = -7307.8181
$ws.Range("H141").Value = 1664.0714
$ws.Range("I141").Value = 1022.8461
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 3068.5383
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = 2111.4617
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22048.607
$ws.Range("J32").Value = 3538.2727
$ws.Range("L32").Value = 3538.2727
$ws.Range("N32").Value = -4112.2727
$ws.Range("H45").Value = 698
$ws.Range("I45").Value = 747.5
$ws.Range("J45").Value = 500
$ws.Range("K45").Value = 747.5
$ws.Range("L45").Value = 500
$ws.Range("M45").Value = -370.5
$ws.Range("N45").Value = -1254
$ws.Range("H74").Value = 1028.75
$ws.Range("I74").Value = 682.5
$ws.Range("J74").Value = 1375
$ws.Range("K74").Value = 682.5
$ws.Range("L74").Value = 1375
$ws.Range("M74").Value = 191.5
$ws.Range("N74").Value = -3123
$ws.Range("H77").Value = 1028.75
$ws.Range("I77").Value = 682.5
$ws.Range("J77").Value = 1375
$ws.Range("K77").Value = 3412.5
$ws.Range("L77").Value = 6875
$ws.Range("M77").Value = 955.5
$ws.Range("N77").Value = -15611
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("H97").Value = 898.3077
$ws.Range("I97").Value = 898.3077
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 898.3077
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -402.3077
$ws.Range("M97").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 517
$ws.Range("I94").Value = 517
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 517
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = -66
$ws.Range("M94").ClearContents()
$ws.Range("H134").Value = 19527.127
$ws.Range("I134").Value = 23333.31
$ws.Range("J134").Value = 2399.3
$ws.Range("K134").Value = 69999.93000000001
$ws.Range("L134").Value = 7197.900000000001
$ws.Range("M134").Value = -67464.93000000001
$ws.Range("N134").Value = -12267.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4653354.5
$ws.Range("I31").Value = 2543.7273
$ws.Range("J31").Value = 9525633
$ws.Range("K31").Value = 2543.7273
$ws.Range("L31").Value = 9525633
$ws.Range("M31").Value = -2248.7273
$ws.Range("N31").Value = -9526223
$ws.Range("H34").Value = 4653354.5
$ws.Range("I34").Value = 2543.7273
$ws.Range("J34").Value = 9525633
$ws.Range("K34").Value = 2543.7273
$ws.Range("L34").Value = 9525633
$ws.Range("M34").Value = -2341.7273
$ws.Range("N34").Value = -9526037
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("H58").Value = 751.4146
$ws.Range("I58").Value = 778.4516
$ws.Range("J58").Value = 667.6
$ws.Range("K58").Value = 778.4516
$ws.Range("L58").Value = 667.6
$ws.Range("M58").Value = -575.4516
$ws.Range("N58").Value = -1073.6
$ws.Range("H132").Value = 2860.8572
$ws.Range("I132").Value = 2286
$ws.Range("J132").Value = 4585.4287
$ws.Range("K132").Value = 6858
$ws.Range("L132").Value = 13756.2861
$ws.Range("M132").Value = -4328
$ws.Range("N132").Value = -18816.2861
$ws.Range("H136").Value = 751.4146
$ws.Range("I136").Value = 778.4516
$ws.Range("J136").Value = 667.6
$ws.Range("K136").Value = 2335.3548
$ws.Range("L136").Value = 2002.8
$ws.Range("M136").Value = 214.6451999999999
$ws.Range("N136").Value = -7102.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1340.0209
$ws.Range("I68").Value = 1204.3914
$ws.Range("J68").Value = 1464.8
$ws.Range("K68").Value = 3613.1742
$ws.Range("L68").Value = 4394.4
$ws.Range("M68").Value = -2802.1742
$ws.Range("N68").Value = -6016.4
$ws.Range("H71").Value = 1340.0209
$ws.Range("I71").Value = 1204.3914
$ws.Range("J71").Value = 1464.8
$ws.Range("K71").Value = 10839.5226
$ws.Range("L71").Value = 13183.2
$ws.Range("M71").Value = -6783.5226
$ws.Range("N71").Value = -21295.2
$ws.Range("H131").Value = 4788.269
$ws.Range("J131").Value = 1133.75
$ws.Range("L131").Value = 3401.25
$ws.Range("N131").Value = -13481.25
$ws.Range("H140").Value = 1155.625
$ws.Range("I140").Value = 608.46155
$ws.Range("J140").Value = 3526.6667
$ws.Range("K140").Value = 1825.38465
$ws.Range("L140").Value = 10580.0001
$ws.Range("M140").Value = 3354.61535
$ws.Range("N140").Value = -20940.0001
$ws.Range("H141").Value = 2079.923
$ws.Range("I141").Value = 1860.8572
$ws.Range("K141").Value = 5582.571599999999
$ws.Range("M141").Value = -402.5715999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6779.5713
$ws.Range("I122").Value = 56300
$ws.Range("J122").Value = 2970.3076
$ws.Range("K122").Value = 168900
$ws.Range("L122").Value = 8910.9228
$ws.Range("M122").Value = -166450
$ws.Range("N122").Value = -13810.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6207.625
$ws.Range("I136").Value = 8793.071
$ws.Range("J136").Value = 2588
$ws.Range("K136").Value = 26379.213
$ws.Range("L136").Value = 7764
$ws.Range("M136").Value = -23829.213
$ws.Range("N136").Value = -12864

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 21873
$ws.Range("I136").Value = 42486
$ws.Range("J136").Value = 1260
$ws.Range("K136").Value = 127458
$ws.Range("L136").Value = 3780
$ws.Range("M136").Value = -124908
$ws.Range("N136").Value = -8880
